$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-6 (columns D through AJ) with the restated (corrected) figures ---

$data = @{
    2 = @{ D=4754; E=124; F=124; G=87; H=67; I=67; J=0; K=3644; L=2651; M=993; N=987; O=6; P=200; Q=75; R=-149; S=91; T=142; U=-67; V=1412; W=2.61; X=1.41; Y=7.03; Z=1.84; AA=266.88; AB=397.52; AC=167; AD=5.82; AE=2539; AF=0.38; AG=15; AH=1.54; AI=8.96; AJ=36700000 }
    3 = @{ D=4824; E=131; F=131; G=89; H=40; I=39; J=0; K=3701; L=2695; M=1007; N=1000; O=6; P=200; Q=140; R=-191; S=68; T=195; U=-55; V=1486; W=2.72; X=0.82; Y=3.93; Z=1.08; AA=267.74; AB=407.77; AC=98; AD=17.21; AE=2572; AF=0.65; AG=15; AH=0.89; AI=15.36; AJ=36700000 }
    4 = @{ D=4770; E=156; F=156; G=108; H=96; I=96; J=0; K=3731; L=2652; M=1078; N=1072; O=6; P=200; Q=181; R=-122; S=-23; T=123; U=58; V=1468; W=3.28; X=2.01; Y=9.3; Z=2.59; AA=245.95; AB=450.98; AC=241; AD=9.51; AE=2757; AF=0.83; AG=20; AH=0.87; AI=8.24; AJ=36700000 }
    5 = @{ D=4692; E=100; F=100; G=-78; H=-42; I=-41; J=-1; K=3771; L=2732; M=1038; N=1032; O=6; P=200; Q=162; R=-224; S=109; T=183; U=-20; V=1570; W=2.13; X=-0.89; Y=-3.89; Z=-1.11; AA=263.19; AB=440.35; AC=-102; AD=-16.63; AE=2598; AF=0.65; AG=5; AH=0.29; AI=-5.26; AJ=36700000 }
    6 = @{ D=4492; E=72; F=72; G=-11; H=-54; I=-54; K=3963; L=2985; M=978; N=971; P=200; Q=12; R=-152; S=125; T=168; U=-156; V=1779; W=1.6; X=-1.21; Y=-5.35; Z=-1.41; AA=305.18; AB=411.35; AC=-134; AD=-10.07; AE=2446; AF=0.55; AG=5; AH=0.37; AI=-4.01; AJ=36700000 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# --- Rows 7-9: the detailed financial figures were erroneous and are removed,        ---
# --- leaving only the index (A), ticker (B) and company name (C) columns populated.  ---
$ws.Range("D7:AJ9").ClearContents()
